$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New daily GSC export rows to append below the existing data (rows 2-71
# already hold 2025-10-07 .. 2025-12-15). We add 2025-12-16 .. 2025-12-22.
$dates  = @("2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22")
$pages  = @(31, 31, 31, 31, 32, 32, 32)

$startRow = 72
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Use Formula with a leading apostrophe so the date-like text is stored
    # as a literal string (matching the workbook's existing convention of
    # shared-string dates) instead of being auto-converted to a date serial
    # number. ClearFormats() afterwards drops the transient quote-prefix
    # style so the cell keeps the sheet's default (unstyled) look.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = "'" + $dates[$i]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = $pages[$i]
}
